# Auto update Excel log
# Append new Bedroom Door proximity events to the "Proximity" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proximity")

# New rows to append (Date, Timestamp, Hour, Location, Value, Status)
$newRows = @(
    @("2026-02-01", "15:08:31", "15:00", "Bedroom Door", "EXIT",  "User EXITED Bedroom"),
    @("2026-02-01", "15:08:36", "15:00", "Bedroom Door", "ENTER", "User ENTERED Bedroom"),
    @("2026-02-01", "15:08:38", "15:00", "Bedroom Door", "EXIT",  "User EXITED Bedroom"),
    @("2026-02-01", "15:09:01", "15:00", "Bedroom Door", "ENTER", "User ENTERED Bedroom"),
    @("2026-02-01", "15:09:03", "15:00", "Bedroom Door", "EXIT",  "User EXITED Bedroom"),
    @("2026-02-01", "15:09:07", "15:00", "Bedroom Door", "ENTER", "User ENTERED Bedroom"),
    @("2026-02-01", "15:09:15", "15:00", "Bedroom Door", "EXIT",  "User EXITED Bedroom"),
    @("2026-02-01", "15:09:17", "15:00", "Bedroom Door", "ENTER", "User ENTERED Bedroom"),
    @("2026-02-01", "15:09:27", "15:00", "Bedroom Door", "EXIT",  "User EXITED Bedroom")
)

$startRow = 9
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Column A holds a date-like string ("2026-02-01"). Assigning it directly
    # would let Excel auto-convert it to a date serial number, so force the
    # cell to text first, then clear the leftover number formatting so the
    # cell keeps using the sheet's default style (matching the other rows).
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rowData[0]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}
